# Auto-generated PowerShell COM-interop edit script
# Daily attendance processing - 2025-11-19 13:35:20

$wb = $excel.ActiveWorkbook
$wsAtt = $wb.Worksheets.Item("Attendance")
$wsSum = $wb.Worksheets.Item("Summary")

# --- 1. Append 19 new attendance rows (576-594) to the Attendance sheet ---
$newDataRange = $wsAtt.Range("A576:K594")
$newDataRange.NumberFormat = "@"

$wsAtt.Cells.Item(576, 1).Value = '221539'
$wsAtt.Cells.Item(576, 2).Value = 'تحريم شوكات مالك'
$wsAtt.Cells.Item(576, 3).Value = 'Year 2'
$wsAtt.Cells.Item(576, 4).Value = 'C1'
$wsAtt.Cells.Item(576, 5).Value = '221539@med.asu.edu.eg'
$wsAtt.Cells.Item(576, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(576, 7).Value = '1'
$wsAtt.Cells.Item(576, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(576, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(576, 10).Value = '11:28:04'
$wsAtt.Cells.Item(576, 11).Value = 'C1'

$wsAtt.Cells.Item(577, 1).Value = '221755'
$wsAtt.Cells.Item(577, 2).Value = 'سعدية عاشق'
$wsAtt.Cells.Item(577, 3).Value = 'Year 2'
$wsAtt.Cells.Item(577, 4).Value = 'C1'
$wsAtt.Cells.Item(577, 5).Value = '221755@med.asu.edu.eg'
$wsAtt.Cells.Item(577, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(577, 7).Value = '1'
$wsAtt.Cells.Item(577, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(577, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(577, 10).Value = '11:28:16'
$wsAtt.Cells.Item(577, 11).Value = 'C1'

$wsAtt.Cells.Item(578, 1).Value = '221833'
$wsAtt.Cells.Item(578, 2).Value = 'صفا محمود صايل صايل'
$wsAtt.Cells.Item(578, 3).Value = 'Year 2'
$wsAtt.Cells.Item(578, 4).Value = 'C1'
$wsAtt.Cells.Item(578, 5).Value = '221833@med.asu.edu.eg'
$wsAtt.Cells.Item(578, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(578, 7).Value = '1'
$wsAtt.Cells.Item(578, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(578, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(578, 10).Value = '11:28:49'
$wsAtt.Cells.Item(578, 11).Value = 'C1'

$wsAtt.Cells.Item(579, 1).Value = '222058'
$wsAtt.Cells.Item(579, 2).Value = 'رفا السيد قسم الله السيد'
$wsAtt.Cells.Item(579, 3).Value = 'Year 2'
$wsAtt.Cells.Item(579, 4).Value = 'C1'
$wsAtt.Cells.Item(579, 5).Value = '222058@med.asu.edu.eg'
$wsAtt.Cells.Item(579, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(579, 7).Value = '1'
$wsAtt.Cells.Item(579, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(579, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(579, 10).Value = '11:31:08'
$wsAtt.Cells.Item(579, 11).Value = 'C1'

$wsAtt.Cells.Item(580, 1).Value = '221904'
$wsAtt.Cells.Item(580, 2).Value = 'عائشه نور شيهو'
$wsAtt.Cells.Item(580, 3).Value = 'Year 2'
$wsAtt.Cells.Item(580, 4).Value = 'C1'
$wsAtt.Cells.Item(580, 5).Value = '221904@med.asu.edu.eg'
$wsAtt.Cells.Item(580, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(580, 7).Value = '1'
$wsAtt.Cells.Item(580, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(580, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(580, 10).Value = '11:31:54'
$wsAtt.Cells.Item(580, 11).Value = 'C1'

$wsAtt.Cells.Item(581, 1).Value = '221948'
$wsAtt.Cells.Item(581, 2).Value = 'سانتينو اتيم شول دينق'
$wsAtt.Cells.Item(581, 3).Value = 'Year 2'
$wsAtt.Cells.Item(581, 4).Value = 'C1'
$wsAtt.Cells.Item(581, 5).Value = '221948@med.asu.edu.eg'
$wsAtt.Cells.Item(581, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(581, 7).Value = '1'
$wsAtt.Cells.Item(581, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(581, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(581, 10).Value = '11:32:33'
$wsAtt.Cells.Item(581, 11).Value = 'C1'

$wsAtt.Cells.Item(582, 1).Value = '221714'
$wsAtt.Cells.Item(582, 2).Value = 'زينب عبد اللطيف بيبى فاروق'
$wsAtt.Cells.Item(582, 3).Value = 'Year 2'
$wsAtt.Cells.Item(582, 4).Value = 'C1'
$wsAtt.Cells.Item(582, 5).Value = '221714@med.asu.edu.eg'
$wsAtt.Cells.Item(582, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(582, 7).Value = '1'
$wsAtt.Cells.Item(582, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(582, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(582, 10).Value = '11:34:55'
$wsAtt.Cells.Item(582, 11).Value = 'C1'

$wsAtt.Cells.Item(583, 1).Value = '221822'
$wsAtt.Cells.Item(583, 2).Value = 'سعادة يوسف عليو'
$wsAtt.Cells.Item(583, 3).Value = 'Year 2'
$wsAtt.Cells.Item(583, 4).Value = 'C1'
$wsAtt.Cells.Item(583, 5).Value = '221822@med.asu.edu.eg'
$wsAtt.Cells.Item(583, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(583, 7).Value = '1'
$wsAtt.Cells.Item(583, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(583, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(583, 10).Value = '11:35:26'
$wsAtt.Cells.Item(583, 11).Value = 'C1'

$wsAtt.Cells.Item(584, 1).Value = '211620'
$wsAtt.Cells.Item(584, 2).Value = 'محمودول اسلام'
$wsAtt.Cells.Item(584, 3).Value = 'Year 2'
$wsAtt.Cells.Item(584, 4).Value = 'C1'
$wsAtt.Cells.Item(584, 5).Value = '211620@med.asu.edu.eg'
$wsAtt.Cells.Item(584, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(584, 7).Value = '1'
$wsAtt.Cells.Item(584, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(584, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(584, 10).Value = '11:35:49'
$wsAtt.Cells.Item(584, 11).Value = 'C1'

$wsAtt.Cells.Item(585, 1).Value = '222076'
$wsAtt.Cells.Item(585, 2).Value = 'ابرار عبد الماجد عبد العزيز عثمان'
$wsAtt.Cells.Item(585, 3).Value = 'Year 2'
$wsAtt.Cells.Item(585, 4).Value = 'C1'
$wsAtt.Cells.Item(585, 5).Value = '222076@med.asu.edu.eg'
$wsAtt.Cells.Item(585, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(585, 7).Value = '1'
$wsAtt.Cells.Item(585, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(585, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(585, 10).Value = '11:36:26'
$wsAtt.Cells.Item(585, 11).Value = 'C1'

$wsAtt.Cells.Item(586, 1).Value = '212442'
$wsAtt.Cells.Item(586, 2).Value = 'رميساء محى الدين الامين الطيب'
$wsAtt.Cells.Item(586, 3).Value = 'Year 2'
$wsAtt.Cells.Item(586, 4).Value = 'C1'
$wsAtt.Cells.Item(586, 5).Value = '212442@med.asu.edu.eg'
$wsAtt.Cells.Item(586, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(586, 7).Value = '1'
$wsAtt.Cells.Item(586, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(586, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(586, 10).Value = '11:36:42'
$wsAtt.Cells.Item(586, 11).Value = 'C1'

$wsAtt.Cells.Item(587, 1).Value = '220304'
$wsAtt.Cells.Item(587, 2).Value = 'احمد الكامل محمد عبدون عثمان'
$wsAtt.Cells.Item(587, 3).Value = 'Year 2'
$wsAtt.Cells.Item(587, 4).Value = 'C1'
$wsAtt.Cells.Item(587, 5).Value = '220304@med.asu.edu.eg'
$wsAtt.Cells.Item(587, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(587, 7).Value = '1'
$wsAtt.Cells.Item(587, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(587, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(587, 10).Value = '11:37:00'
$wsAtt.Cells.Item(587, 11).Value = 'C1'

$wsAtt.Cells.Item(588, 1).Value = '220967'
$wsAtt.Cells.Item(588, 2).Value = 'لارا حربي عبدالله الزيادات'
$wsAtt.Cells.Item(588, 3).Value = 'Year 2'
$wsAtt.Cells.Item(588, 4).Value = 'C1'
$wsAtt.Cells.Item(588, 5).Value = '220967@med.asu.edu.eg'
$wsAtt.Cells.Item(588, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(588, 7).Value = '1'
$wsAtt.Cells.Item(588, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(588, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(588, 10).Value = '11:37:18'
$wsAtt.Cells.Item(588, 11).Value = 'C1'

$wsAtt.Cells.Item(589, 1).Value = '212543'
$wsAtt.Cells.Item(589, 2).Value = 'زينب سيف الدين محمد ادم'
$wsAtt.Cells.Item(589, 3).Value = 'Year 2'
$wsAtt.Cells.Item(589, 4).Value = 'C1'
$wsAtt.Cells.Item(589, 5).Value = '212543@med.asu.edu.eg'
$wsAtt.Cells.Item(589, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(589, 7).Value = '1'
$wsAtt.Cells.Item(589, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(589, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(589, 10).Value = '11:37:37'
$wsAtt.Cells.Item(589, 11).Value = 'C1'

$wsAtt.Cells.Item(590, 1).Value = '222028'
$wsAtt.Cells.Item(590, 2).Value = 'هاجر عبد الحفيظ سيد صالح'
$wsAtt.Cells.Item(590, 3).Value = 'Year 2'
$wsAtt.Cells.Item(590, 4).Value = 'C1'
$wsAtt.Cells.Item(590, 5).Value = '222028@med.asu.edu.eg'
$wsAtt.Cells.Item(590, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(590, 7).Value = '1'
$wsAtt.Cells.Item(590, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(590, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(590, 10).Value = '11:37:49'
$wsAtt.Cells.Item(590, 11).Value = 'C1'

$wsAtt.Cells.Item(591, 1).Value = '222003'
$wsAtt.Cells.Item(591, 2).Value = 'اسراء بدر الدين جعفر عثمان'
$wsAtt.Cells.Item(591, 3).Value = 'Year 2'
$wsAtt.Cells.Item(591, 4).Value = 'C1'
$wsAtt.Cells.Item(591, 5).Value = '222003@med.asu.edu.eg'
$wsAtt.Cells.Item(591, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(591, 7).Value = '1'
$wsAtt.Cells.Item(591, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(591, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(591, 10).Value = '11:37:54'
$wsAtt.Cells.Item(591, 11).Value = 'C1'

$wsAtt.Cells.Item(592, 1).Value = '221000'
$wsAtt.Cells.Item(592, 2).Value = 'ابوبكر محمد قايد الثوابي'
$wsAtt.Cells.Item(592, 3).Value = 'Year 2'
$wsAtt.Cells.Item(592, 4).Value = 'C1'
$wsAtt.Cells.Item(592, 5).Value = '221000@med.asu.edu.eg'
$wsAtt.Cells.Item(592, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(592, 7).Value = '1'
$wsAtt.Cells.Item(592, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(592, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(592, 10).Value = '11:38:06'
$wsAtt.Cells.Item(592, 11).Value = 'C1'

$wsAtt.Cells.Item(593, 1).Value = '222004'
$wsAtt.Cells.Item(593, 2).Value = 'احمد ايمن احمد بشير'
$wsAtt.Cells.Item(593, 3).Value = 'Year 2'
$wsAtt.Cells.Item(593, 4).Value = 'C1'
$wsAtt.Cells.Item(593, 5).Value = '222004@med.asu.edu.eg'
$wsAtt.Cells.Item(593, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(593, 7).Value = '1'
$wsAtt.Cells.Item(593, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(593, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(593, 10).Value = '11:38:18'
$wsAtt.Cells.Item(593, 11).Value = 'C1'

$wsAtt.Cells.Item(594, 1).Value = '222053'
$wsAtt.Cells.Item(594, 2).Value = 'صباح سيف الدين عثمان اسحق'
$wsAtt.Cells.Item(594, 3).Value = 'Year 2'
$wsAtt.Cells.Item(594, 4).Value = 'C1'
$wsAtt.Cells.Item(594, 5).Value = '222053@med.asu.edu.eg'
$wsAtt.Cells.Item(594, 6).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(594, 7).Value = '1'
$wsAtt.Cells.Item(594, 8).Value = 'BIOCHEMISTRY LAB/CBL'
$wsAtt.Cells.Item(594, 9).Value = '19/11/2025'
$wsAtt.Cells.Item(594, 10).Value = '11:38:41'
$wsAtt.Cells.Item(594, 11).Value = 'C1'

$newDataRange.ClearFormats()

# --- 2. Update AutoFilter range and dimension on Attendance sheet ---
$wsAtt.AutoFilterMode = $false
$wsAtt.Range("A1:K594").AutoFilter()

# --- 3. Update the _FilterDatabase defined name for Attendance sheet ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "=Attendance!`$A`$1:`$K`$594"
    }
}

# --- 4. Widen columns F and H on Attendance sheet (14 -> 22) ---
$wsAtt.Columns.Item(6).ColumnWidth = 21.1666666666667
$wsAtt.Columns.Item(8).ColumnWidth = 21.1666666666667

# --- 5. Update computed Summary-sheet stats for the 19 affected students ---
# Row 37
$wsSum.Cells.Item(37, 7).Value = '6.9%'
$wsSum.Cells.Item(37, 9).Value = 21
$wsSum.Cells.Item(37, 14).Value = 2
$wsSum.Cells.Item(37, 15).Value = 9
$wsSum.Cells.Item(37, 24).Value = 1
$wsSum.Cells.Item(37, 25).Value = 1

# Row 71
$wsSum.Cells.Item(71, 6).Value = 'High Risk'
$wsSum.Cells.Item(71, 6).Interior.Color = 8158463
$wsSum.Cells.Item(71, 7).Value = '17.2%'
$wsSum.Cells.Item(71, 9).Value = 18
$wsSum.Cells.Item(71, 14).Value = 5
$wsSum.Cells.Item(71, 15).Value = 6
$wsSum.Cells.Item(71, 24).Value = 1
$wsSum.Cells.Item(71, 25).Value = 1

# Row 74
$wsSum.Cells.Item(74, 6).Value = 'High Risk'
$wsSum.Cells.Item(74, 6).Interior.Color = 8158463
$wsSum.Cells.Item(74, 7).Value = '17.2%'
$wsSum.Cells.Item(74, 9).Value = 18
$wsSum.Cells.Item(74, 14).Value = 5
$wsSum.Cells.Item(74, 15).Value = 6
$wsSum.Cells.Item(74, 24).Value = 1
$wsSum.Cells.Item(74, 25).Value = 1

# Row 75
$wsSum.Cells.Item(75, 6).Value = 'High Risk'
$wsSum.Cells.Item(75, 6).Interior.Color = 8158463
$wsSum.Cells.Item(75, 7).Value = '17.2%'
$wsSum.Cells.Item(75, 9).Value = 18
$wsSum.Cells.Item(75, 14).Value = 5
$wsSum.Cells.Item(75, 15).Value = 6
$wsSum.Cells.Item(75, 24).Value = 1
$wsSum.Cells.Item(75, 25).Value = 1

# Row 93
$wsSum.Cells.Item(93, 7).Value = '13.8%'
$wsSum.Cells.Item(93, 9).Value = 19
$wsSum.Cells.Item(93, 14).Value = 4
$wsSum.Cells.Item(93, 15).Value = 7
$wsSum.Cells.Item(93, 24).Value = 1
$wsSum.Cells.Item(93, 25).Value = 1

# Row 96
$wsSum.Cells.Item(96, 7).Value = '13.8%'
$wsSum.Cells.Item(96, 9).Value = 19
$wsSum.Cells.Item(96, 14).Value = 4
$wsSum.Cells.Item(96, 15).Value = 7
$wsSum.Cells.Item(96, 24).Value = 1
$wsSum.Cells.Item(96, 25).Value = 1

# Row 147
$wsSum.Cells.Item(147, 7).Value = '20.7%'
$wsSum.Cells.Item(147, 9).Value = 17
$wsSum.Cells.Item(147, 14).Value = 6
$wsSum.Cells.Item(147, 15).Value = 5
$wsSum.Cells.Item(147, 24).Value = 1
$wsSum.Cells.Item(147, 25).Value = 1

# Row 180
$wsSum.Cells.Item(180, 7).Value = '3.4%'
$wsSum.Cells.Item(180, 9).Value = 22
$wsSum.Cells.Item(180, 14).Value = 1
$wsSum.Cells.Item(180, 15).Value = 10
$wsSum.Cells.Item(180, 24).Value = 1
$wsSum.Cells.Item(180, 25).Value = 1

# Row 184
$wsSum.Cells.Item(184, 7).Value = '20.7%'
$wsSum.Cells.Item(184, 9).Value = 17
$wsSum.Cells.Item(184, 14).Value = 6
$wsSum.Cells.Item(184, 15).Value = 5
$wsSum.Cells.Item(184, 24).Value = 1
$wsSum.Cells.Item(184, 25).Value = 1

# Row 195
$wsSum.Cells.Item(195, 6).Value = 'High Risk'
$wsSum.Cells.Item(195, 6).Interior.Color = 8158463
$wsSum.Cells.Item(195, 7).Value = '17.2%'
$wsSum.Cells.Item(195, 9).Value = 18
$wsSum.Cells.Item(195, 14).Value = 5
$wsSum.Cells.Item(195, 15).Value = 6
$wsSum.Cells.Item(195, 24).Value = 1
$wsSum.Cells.Item(195, 25).Value = 1

# Row 197
$wsSum.Cells.Item(197, 7).Value = '13.8%'
$wsSum.Cells.Item(197, 9).Value = 19
$wsSum.Cells.Item(197, 14).Value = 4
$wsSum.Cells.Item(197, 15).Value = 7
$wsSum.Cells.Item(197, 24).Value = 1
$wsSum.Cells.Item(197, 25).Value = 1

# Row 209
$wsSum.Cells.Item(209, 7).Value = '13.8%'
$wsSum.Cells.Item(209, 9).Value = 19
$wsSum.Cells.Item(209, 14).Value = 4
$wsSum.Cells.Item(209, 15).Value = 7
$wsSum.Cells.Item(209, 24).Value = 1
$wsSum.Cells.Item(209, 25).Value = 1

# Row 220
$wsSum.Cells.Item(220, 7).Value = '10.3%'
$wsSum.Cells.Item(220, 9).Value = 20
$wsSum.Cells.Item(220, 14).Value = 3
$wsSum.Cells.Item(220, 15).Value = 8
$wsSum.Cells.Item(220, 24).Value = 1
$wsSum.Cells.Item(220, 25).Value = 1

# Row 232
$wsSum.Cells.Item(232, 7).Value = '20.7%'
$wsSum.Cells.Item(232, 9).Value = 17
$wsSum.Cells.Item(232, 14).Value = 6
$wsSum.Cells.Item(232, 15).Value = 5
$wsSum.Cells.Item(232, 24).Value = 1
$wsSum.Cells.Item(232, 25).Value = 1

# Row 233
$wsSum.Cells.Item(233, 7).Value = '27.6%'
$wsSum.Cells.Item(233, 9).Value = 15
$wsSum.Cells.Item(233, 14).Value = 8
$wsSum.Cells.Item(233, 15).Value = 3
$wsSum.Cells.Item(233, 24).Value = 1
$wsSum.Cells.Item(233, 25).Value = 1

# Row 238
$wsSum.Cells.Item(238, 7).Value = '20.7%'
$wsSum.Cells.Item(238, 9).Value = 17
$wsSum.Cells.Item(238, 14).Value = 6
$wsSum.Cells.Item(238, 15).Value = 5
$wsSum.Cells.Item(238, 24).Value = 1
$wsSum.Cells.Item(238, 25).Value = 1

# Row 243
$wsSum.Cells.Item(243, 7).Value = '20.7%'
$wsSum.Cells.Item(243, 9).Value = 17
$wsSum.Cells.Item(243, 14).Value = 6
$wsSum.Cells.Item(243, 15).Value = 5
$wsSum.Cells.Item(243, 24).Value = 1
$wsSum.Cells.Item(243, 25).Value = 1

# Row 245
$wsSum.Cells.Item(245, 7).Value = '20.7%'
$wsSum.Cells.Item(245, 9).Value = 17
$wsSum.Cells.Item(245, 14).Value = 6
$wsSum.Cells.Item(245, 15).Value = 5
$wsSum.Cells.Item(245, 24).Value = 1
$wsSum.Cells.Item(245, 25).Value = 1

# Row 248
$wsSum.Cells.Item(248, 7).Value = '10.3%'
$wsSum.Cells.Item(248, 9).Value = 20
$wsSum.Cells.Item(248, 14).Value = 3
$wsSum.Cells.Item(248, 15).Value = 8
$wsSum.Cells.Item(248, 24).Value = 1
$wsSum.Cells.Item(248, 25).Value = 1
